# Applies the Oct 22 2023 GitHub Actions "cryptos list" refresh:
# updated Price (D) / Volume(1h) (E) figures for most rows, and a
# reordering of the Kaspa / WEMIXToken rows (43 <-> 44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.886.32'
$ws.Range('E2').Value = '  +1.20%  '

# Row 3
$ws.Range('D3').Value = '1.626.43'
$ws.Range('E3').Value = '  +1.83%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').Value = '''214.16'
$ws.Range('E5').Value = '  +1.04%  '

# Row 6
$ws.Range('D6').Value = '''0.519'
$ws.Range('E6').Value = '  +1.11%  '

# Row 7
$ws.Range('E7').Value = '  -0.10%  '

# Row 8
$ws.Range('E8').Value = '  +11.08%  '

# Row 9
$ws.Range('E9').Value = '  +3.43%  '

# Row 10
$ws.Range('E10').Value = '  +2.52%  '

# Row 11
$ws.Range('D11').Value = '''0.0917'
$ws.Range('E11').Value = '  +0.90%  '

# Row 12
$ws.Range('D12').Value = '1.859.53'
$ws.Range('E12').Value = '  +1.88%  '

# Row 13
$ws.Range('D13').Value = '1.627.27'
$ws.Range('E13').Value = '  +1.91%  '

# Row 14
$ws.Range('D14').Value = '''0.570'
$ws.Range('E14').Value = '  +6.37%  '

# Row 15
$ws.Range('D15').Value = '''3.89'
$ws.Range('E15').Value = '  +4.44%  '

# Row 16
$ws.Range('D16').Value = '29.919.04'
$ws.Range('E16').Value = '  +1.29%  '

# Row 17
$ws.Range('D17').Value = '''9.06'
$ws.Range('E17').Value = '  +19.96%  '

# Row 18
$ws.Range('D18').Value = '''64.93'
$ws.Range('E18').Value = '  +1.84%  '

# Row 19
$ws.Range('D19').Value = '''246.46'
$ws.Range('E19').Value = '  +2.49%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0704'
$ws.Range('E20').Value = '  +1.70%  '

# Row 21
$ws.Range('E21').Value = '  -0.07%  '

# Row 22
$ws.Range('D22').Value = '''4.13'
$ws.Range('E22').Value = '  +3.57%  '

# Row 23
$ws.Range('D23').Value = '''9.60'
$ws.Range('E23').Value = '  +4.06%  '

# Row 24
$ws.Range('E24').Value = '  +0.59%  '

# Row 25
$ws.Range('D25').Value = '''158.24'
$ws.Range('E25').Value = '  +2.12%  '

# Row 26
$ws.Range('D26').Value = '''15.69'
$ws.Range('E26').Value = '  +2.40%  '

# Row 27
$ws.Range('E27').Value = '  +2.44%  '

# Row 28
$ws.Range('E28').Value = '  +3.05%  '

# Row 29
$ws.Range('E29').Value = '  -0.03%  '

# Row 30
$ws.Range('E30').Value = '  +2.72%  '

# Row 31
$ws.Range('E31').Value = '  +6.00%  '

# Row 32
$ws.Range('E32').Value = '  +4.14%  '

# Row 33
$ws.Range('E33').Value = '  +1.88%  '

# Row 34
$ws.Range('D34').Value = '1.429.10'
$ws.Range('E34').Value = '  +0.02%  '

# Row 35
$ws.Range('E35').Value = '  +6.87%  '

# Row 36
$ws.Range('E36').Value = '  +1.10%  '

# Row 37
$ws.Range('D37').Value = '''2.88'
$ws.Range('E37').Value = '  +1.90%  '

# Row 38
$ws.Range('E38').Value = '  -0.40%  '

# Row 39
$ws.Range('E39').Value = '  +3.27%  '

# Row 40
$ws.Range('D40').Value = '''0.555'
$ws.Range('E40').Value = '  +3.19%  '

# Row 41
$ws.Range('D41').Value = '''0.831'
$ws.Range('E41').Value = '  +3.92%  '

# Row 42
$ws.Range('D42').Value = '''55.16'
$ws.Range('E42').Value = '  +3.37%  '

# Row 43
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '''1.06'
$ws.Range('E43').Value = '  +7.58%  '

# Row 44
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '''0.0496'
$ws.Range('E44').Value = '  +1.03%  '

# Row 45
$ws.Range('D45').Value = '''1.97'
$ws.Range('E45').Value = '  +0.23%  '

# Row 46
$ws.Range('D46').Value = '''70.21'
$ws.Range('E46').Value = '  +7.18%  '

# Row 47
$ws.Range('E47').Value = '  -0.06%  '

# Row 49
$ws.Range('D49').Value = '1.767.26'
$ws.Range('E49').Value = '  +1.73%  '

# Row 50
$ws.Range('E50').Value = '  +4.08%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0107'
$ws.Range('E51').Value = '  +1.33%  '
